$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.025.61'

$ws.Range('D3').Value = '2.602.35'
$ws.Range('E3').Value = '  -0.46%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.30'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.06'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -2.48%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.545'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -0.66%  '

$ws.Range('D9').Value = '2.598.61'
$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('E10').Value = '  +1.20%  '

$ws.Range('E11').Value = '  -0.02%  '

$ws.Range('E12').Value = '  -1.62%  '

$ws.Range('E13').Value = '  -3.43%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.11'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -2.83%  '

$ws.Range('D15').Value = '3.074.52'
$ws.Range('E15').Value = '  -0.49%  '

$ws.Range('E16').Value = '  -3.12%  '

$ws.Range('D17').Value = '66.900.16'
$ws.Range('E17').Value = '  -0.99%  '

$ws.Range('D18').Value = '2.595.85'
$ws.Range('E18').Value = '  -0.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '367.56'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.40%  '

$ws.Range('E20').Value = '  -2.18%  '

$ws.Range('E21').Value = '  -3.29%  '

$ws.Range('E22').Value = '  -0.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.82'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.93%  '

$ws.Range('E24').Value = '  -3.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.48'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +5.07%  '

$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -0.87%  '

$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +0.16%  '

$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '581.71'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +0.08%  '

$ws.Range('D31').Value = '0.0₃0984'
$ws.Range('E31').Value = '  -6.79%  '

$ws.Range('E32').Value = '  -5.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.65'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -3.52%  '

$ws.Range('E34').Value = '  -3.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.05%  '

$ws.Range('E36').Value = '  -4.22%  '

$ws.Range('E37').Value = '  -3.02%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.62'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +0.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.98'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -2.30%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.364'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -1.75%  '

$ws.Range('E41').Value = '  -0.69%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.22'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -3.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.56'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -3.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.10'
$ws.Range('D44').NumberFormat = 'General'

$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '152.72'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -2.57%  '

$ws.Range('D47').Value = '0.0₆0284'
$ws.Range('E47').Value = '  -2.23%  '

$ws.Range('E48').Value = '  -1.48%  '

$ws.Range('E49').Value = '  -1.65%  '

$ws.Range('E50').Value = '  -4.01%  '

$ws.Range('E51').Value = '  +1.49%  '
